$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set updated/added numeric values
$ws.Range("D2").Value = -0.462
$ws.Range("E2").Value = 0.234
$ws.Range("G2").Value = -1.300567107750473
$ws.Range("H2").Value = -1.300567107750473
$ws.Range("I2").Value = -1.321361058601134
$ws.Range("J2").Value = -1.295350014140482
$ws.Range("K2").Value = -0.615
$ws.Range("L2").Value = -1.162570888468809
$ws.Range("M2").Value = 0.198
$ws.Range("N2").Value = 0.01554160125588697
$ws.Range("O2").Value = -0.3219512195121951
$ws.Range("P2").Value = 0.198
$ws.Range("Q2").Value = 0.01554160125588697
$ws.Range("R2").Value = -0.3219512195121951
$ws.Range("U2").Value = 2.731
$ws.Range("V2").Value = 0.2143642072213501
$ws.Range("W2").Value = 0.01449404761904762
$ws.Range("X2").Value = 0.05921155817013783
$ws.Range("Y2").Value = -0.04471751055109021
$ws.Range("Z2").Value = 0.02356977365888434
$ws.Range("AA2").Value = -0.002380952380952381
$ws.Range("AB2").Value = 0.05921155817013783
$ws.Range("AC2").Value = -0.06159251055109021
$ws.Range("AD2").Value = 0.409
$ws.Range("AF2").Value = 0.409
$ws.Range("AG2").Value = -2.322
$ws.Range("AH2").Value = 0.03110502699825081
$ws.Range("AI2").Value = 0.01770639421619983
$ws.Range("AJ2").Value = -0.2228834709157228
$ws.Range("AK2").Value = -0.1140023566378633
$ws.Range("AL2").Value = 0.031
$ws.Range("AM2").Value = 0.031
$ws.Range("AN2").Value = -1.025062656641604
$ws.Range("AO2").Value = -22.54838709677419
$ws.Range("AP2").Value = 5.819548872180451
$ws.Range("AQ2").Value = -22.54838709677419
$ws.Range("D3").Value = -0.462
$ws.Range("G3").Value = 0.2913256955810147
$ws.Range("H3").Value = 0.2913256955810147
$ws.Range("I3").Value = 0.2242225859247136
$ws.Range("J3").Value = 0.2242225859247136
$ws.Range("K3").Value = 0.135
$ws.Range("L3").Value = 0.220949263502455
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = -0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("U3").Value = 0.148
$ws.Range("V3").Value = 0.06666666666666665
$ws.Range("W3").Value = 0.03375
$ws.Range("X3").Value = 0.05921155817013783
$ws.Range("Y3").Value = -0.02546155817013783
$ws.Range("Z3").Value = 0.1378921236741142
$ws.Range("AA3").Value = 0.03091852854886031
$ws.Range("AB3").Value = 0.05921155817013783
$ws.Range("AC3").Value = -0.02829302962127753
$ws.Range("AD3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = -0.148
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = -0.07142857142857142
$ws.Range("AK3").Value = -0.03707414829659319
$ws.Range("AL3").Value = 0.001
$ws.Range("AM3").Value = 0.001
$ws.Range("AN3").Value = 0
$ws.Range("AO3").Value = 137
$ws.Range("AP3").Value = -0.6727272727272727
$ws.Range("AQ3").Value = 137
$ws.Range("D4").Value = 0.105
$ws.Range("E4").Value = 0.234
$ws.Range("K4").Value = 0.119
$ws.Range("L4").Value = 0.7041420118343195
$ws.Range("M4").Value = 0.198
$ws.Range("N4").Value = 0.05484764542936289
$ws.Range("O4").Value = 1.663865546218487
$ws.Range("P4").Value = 0.198
$ws.Range("Q4").Value = 0.05484764542936289
$ws.Range("R4").Value = 1.663865546218487
$ws.Range("U4").Value = 1.67
$ws.Range("V4").Value = 0.4626038781163435
$ws.Range("W4").Value = 0.03801916932907348
$ws.Range("X4").Value = 0.05921155817013783
$ws.Range("Y4").Value = -0.02119238884106435
$ws.Range("Z4").Value = 0.1251851851851852
$ws.Range("AB4").Value = 0.05921155817013783
$ws.Range("AC4").Value = -0.05921155817013783
$ws.Range("AG4").Value = -1.67
$ws.Range("AJ4").Value = -0.8608247422680412
$ws.Range("AK4").Value = -1.403361344537815
$ws.Range("D5").Value = -0.512
$ws.Range("G5").Value = -23
$ws.Range("H5").Value = -23
$ws.Range("I5").Value = -7
$ws.Range("J5").Value = -7
$ws.Range("K5").Value = -0.021
$ws.Range("L5").Value = -7
$ws.Range("W5").Value = -0.004761904761904762
$ws.Range("X5").Value = 0.05921155817013783
$ws.Range("Y5").Value = -0.06397346293204259
$ws.Range("Z5").Value = 0.0006802721088435374
$ws.Range("AA5").Value = -0.004761904761904762
$ws.Range("AB5").Value = 0.05921155817013783
$ws.Range("AC5").Value = -0.06397346293204259
$ws.Range("G6").Value = 3.137795275590551
$ws.Range("H6").Value = 3.137795275590551
$ws.Range("I6").Value = 3.208661417322834
$ws.Range("J6").Value = 3.208661417322834
$ws.Range("K6").Value = -0.848
$ws.Range("L6").Value = 3.338582677165354
$ws.Range("O6").Value = 0
$ws.Range("R6").Value = 0
$ws.Range("U6").Value = 0.913
$ws.Range("V6").Value = 0.1663023679417122
$ws.Range("W6").Value = -0.07008264462809917
$ws.Range("X6").Value = 0.06218625067164404
$ws.Range("Y6").Value = -0.1322688952997432
$ws.Range("Z6").Value = -0.02072961723659512
$ws.Range("AA6").Value = -0.06651432302293316
$ws.Range("AB6").Value = 0.06174418378588246
$ws.Range("AC6").Value = -0.1282585068088156
$ws.Range("AD6").Value = 0.409
$ws.Range("AF6").Value = 0.409
$ws.Range("AG6").Value = -0.504
$ws.Range("AH6").Value = 0.06933378538735378
$ws.Range("AI6").Value = 0.0349303954223247
$ws.Range("AJ6").Value = -0.1010830324909747
$ws.Range("AK6").Value = -0.04668395702111893
$ws.Range("AL6").Value = 0.03
$ws.Range("AM6").Value = 0.03
$ws.Range("AN6").Value = -0.6607431340872374
$ws.Range("AO6").Value = -27.16666666666666
$ws.Range("AP6").Value = 0.81421647819063
$ws.Range("AQ6").Value = -27.16666666666666

# Clear cells removed in the diff
$ws.Range("E3").ClearContents()
$ws.Range("T3").ClearContents()
$ws.Range("AN5").ClearContents()
$ws.Range("AP5").ClearContents()
$ws.Range("E6").ClearContents()
